# Refactor invalid path detection and fix ScanDate to session timestamp
#
# - Row 2 (LPRIME/commons-io) and Row 4 (LPRIME/poi) ScanDate (col J) now use
#   a single session timestamp instead of the old per-row processing time.
# - Row 5 (previously NOTPRIME2 / UNC-access-denied) becomes a new LPRIME row
#   describing an invalid/root path scan result.
# - Former row 5's NOTPRIME2 data is preserved, but moved down to row 6 with
#   the updated session ScanDate.
# - Two brand-new LPRIME rows (7 and 8) are appended, covering additional
#   invalid-path cases (no path at all, and an "N\A" path).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: FileModificationDate unchanged, ScanDate -> session timestamp ---
$ws.Cells.Item(2, 10).Value2 = "2025-09-09 21:40:41"

# --- Row 4: ScanDate -> session timestamp ---
$ws.Cells.Item(4, 10).Value2 = "2025-09-09 21:40:42"

# --- Row 5: was NOTPRIME2 / UNC access denied row; now an LPRIME invalid-path row ---
$ws.Cells.Item(5, 1).Value2 = "LPRIME"
$ws.Cells.Item(5, 2).Value2 = "Windows Server 2019"
$ws.Cells.Item(5, 3).Value2 = "D:\Docs\Projects\Code\VulnAnalysisTool\"
$ws.Cells.Item(5, 4).Value2 = "X"
$ws.Cells.Item(5, 5).Value2 = "X"
$ws.Cells.Item(5, 9).ClearContents() | Out-Null
$ws.Cells.Item(5, 10).Value2 = "2025-09-09 21:40:42"

# --- Row 6 (new): the old NOTPRIME2 / UNC-access-denied data, re-added below ---
$ws.Cells.Item(6, 1).Value2 = "NOTPRIME2"
$ws.Cells.Item(6, 2).Value2 = "Windows Server 2019"
$ws.Cells.Item(6, 3).Value2 = "D:\Docs\Projects\Code\VulnAnalysisTool\deps\commons-collections4-4.5.0.jar"
$ws.Cells.Item(6, 4).Value2 = "X"
$ws.Cells.Item(6, 9).Value2 = "UNC access denied - cannot determine file existence"
$ws.Cells.Item(6, 10).Value2 = "2025-09-09 21:40:42"

# --- Row 7 (new): LPRIME invalid-path row, no XTRACT_PATH value ---
$ws.Cells.Item(7, 1).Value2 = "LPRIME"
$ws.Cells.Item(7, 2).Value2 = "Windows Server 2019"
$ws.Cells.Item(7, 4).Value2 = "X"
$ws.Cells.Item(7, 5).Value2 = "X"
$ws.Cells.Item(7, 10).Value2 = "2025-09-09 21:40:47"

# --- Row 8 (new): LPRIME invalid-path row, XTRACT_PATH = N\A ---
$ws.Cells.Item(8, 1).Value2 = "LPRIME"
$ws.Cells.Item(8, 2).Value2 = "Windows Server 2019"
$ws.Cells.Item(8, 3).Value2 = "N\A"
$ws.Cells.Item(8, 4).Value2 = "X"
$ws.Cells.Item(8, 5).Value2 = "X"
$ws.Cells.Item(8, 10).Value2 = "2025-09-09 21:40:47"

# Column F/H widths shrink (content got shorter after the ScanDate/path rework);
# approximate the post-edit AutoFit result.
$ws.Columns("F:F").ColumnWidth = 19.75
$ws.Columns("H:H").ColumnWidth = 16.92

# Selection ends on H5 after the edit.
$ws.Range("H5").Select() | Out-Null
